$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(5)   # "Nationalite" donor sheet for cell styles
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Situation_annee_precedente"

# -- Row formats: copy style-only from the donor sheet rows that already use
#    the exact same style indices (1 = blank row, 2 = header, 3/4 = zebra data) --
$src.Range("A1:D1").Copy() | Out-Null
$ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:D22").PasteSpecial(-4122) | Out-Null

$src.Range("A2:D2").Copy() | Out-Null
$ws.Range("A2:D2").PasteSpecial(-4122) | Out-Null

$src.Range("A3:D3").Copy() | Out-Null
$ws.Range("A3:D3").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:D7").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:D9").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:D11").PasteSpecial(-4122) | Out-Null
$ws.Range("A13:D13").PasteSpecial(-4122) | Out-Null
$ws.Range("A15:D15").PasteSpecial(-4122) | Out-Null
$ws.Range("A17:D17").PasteSpecial(-4122) | Out-Null
$ws.Range("A19:D19").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:D21").PasteSpecial(-4122) | Out-Null

$src.Range("A4:D4").Copy() | Out-Null
$ws.Range("A4:D4").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:D6").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:D8").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:D10").PasteSpecial(-4122) | Out-Null
$ws.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$ws.Range("A14:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null
$ws.Range("A18:D18").PasteSpecial(-4122) | Out-Null
$ws.Range("A20:D20").PasteSpecial(-4122) | Out-Null

# -- Cell values --
$ws.Range("B2").Value = 'Situation année précédente (code)'
$ws.Range("C2").Value = 'Situation année précédente (lib.)'
$ws.Range("B3").Value = 'A'
$ws.Range("C3").Value = 'Enseignement secondaire (y compris par correspondance)'
$ws.Range("B4").Value = 'B'
$ws.Range("C4").Value = 'BTS'
$ws.Range("B5").Value = 'C'
$ws.Range("C5").Value = 'IUT'
$ws.Range("B6").Value = 'D'
$ws.Range("C6").Value = 'CPGE (non inscrit à l''université)'
$ws.Range("B7").Value = 'E'
$ws.Range("C7").Value = 'Ecole d''ingénieur (universitaire ou non)'
$ws.Range("B8").Value = 'F'
$ws.Range("C8").Value = 'IUFM'
$ws.Range("B9").Value = 'G'
$ws.Range("C9").Value = 'Enseignement par correspondance'
$ws.Range("B10").Value = 'H'
$ws.Range("C10").Value = 'Université (hors IUT, IUFM, Ecole d''ingénieur universitaire)'
$ws.Range("B11").Value = 'J'
$ws.Range("C11").Value = 'Ecole de MANAGEMENT (commerce, gestion)'
$ws.Range("B12").Value = 'K'
$ws.Range("C12").Value = 'Autre établissement SISE'
$ws.Range("B13").Value = 'L'
$ws.Range("C13").Value = 'Etablissement (hors université) préparant aux concours paramédicaux'
$ws.Range("B14").Value = 'M'
$ws.Range("C14").Value = 'ESPE'
$ws.Range("B15").Value = 'P'
$ws.Range("C15").Value = 'Etablissement étranger d''enseignement supérieur ou secondaire'
$ws.Range("B16").Value = 'Q'
$ws.Range("C16").Value = 'Etab Etranger ens secondaire'
$ws.Range("B17").Value = 'R'
$ws.Range("C17").Value = 'Etab etranger ens supérieur'
$ws.Range("B18").Value = 'S'
$ws.Range("C18").Value = 'Autre établissement ou cursus (hors secondaire, STS, Ingénieur, CPGE, Universitaire, IUT, IUFM, établissements préparant'
$ws.Range("B19").Value = 'T'
$ws.Range("C19").Value = 'Non scolarisé l''année précédente et jamais entré dans l''enseignement supérieur (prise d''étude différée)'
$ws.Range("B20").Value = 'U'
$ws.Range("C20").Value = 'Non scolarisé l''année précédente mais précédemment entré dans l''enseignement supérieur, universitaire ou non (reprise d'''
$ws.Range("B21").Value = 'V'
$ws.Range("C21").Value = 'Instituts catholiques'

# -- Row heights (points) --
$ws.Rows.Item(1).RowHeight = 14.3991
$ws.Rows.Item(2).RowHeight = 23.9985
$ws.Rows.Item(3).RowHeight = 19.7321
$ws.Rows.Item(4).RowHeight = 19.7321
$ws.Rows.Item(5).RowHeight = 19.7321
$ws.Rows.Item(6).RowHeight = 19.7321
$ws.Rows.Item(7).RowHeight = 19.7321
$ws.Rows.Item(8).RowHeight = 19.7321
$ws.Rows.Item(9).RowHeight = 19.7321
$ws.Rows.Item(10).RowHeight = 19.7321
$ws.Rows.Item(11).RowHeight = 19.7321
$ws.Rows.Item(12).RowHeight = 19.7321
$ws.Rows.Item(13).RowHeight = 19.7321
$ws.Rows.Item(14).RowHeight = 19.7321
$ws.Rows.Item(15).RowHeight = 19.7321
$ws.Rows.Item(16).RowHeight = 19.7321
$ws.Rows.Item(17).RowHeight = 19.7321
$ws.Rows.Item(18).RowHeight = 19.7321
$ws.Rows.Item(19).RowHeight = 19.7321
$ws.Rows.Item(20).RowHeight = 19.7321
$ws.Rows.Item(21).RowHeight = 19.7321
$ws.Rows.Item(22).RowHeight = 28.7982

# -- Column widths (chars); engine snaps to integer pixel grid like real Excel --
$ws.Range("A1").ColumnWidth = 15.1667
$ws.Range("B1:C1").ColumnWidth = 9.8333
$ws.Range("D1").ColumnWidth = 3.8333

Write-Host "Sheet created:" $ws.Name
